$d = $word.ActiveDocument

# Update the date line (unique text in the document header paragraph)
$d.Content.Find.Execute("2024-02-03 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-02-04 Sunday", 2)

# The table contains 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17)
# has visible division problems, the rest are blank answer rows.
# Replace cell contents positionally to avoid collisions between old/new
# values that overlap across different cells (e.g. "68÷2=" is both an old
# value and a new value in different cells).
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)
$values = @(
    @("66÷3=", "32÷9=", "98÷4=", "40÷6=", "95÷2="),
    @("58÷9=", "17÷7=", "95÷6=", "68÷2=", "93÷3="),
    @("71÷6=", "79÷3=", "82÷2=", "49÷7=", "15÷7="),
    @("53÷2=", "40÷8=", "32÷9=", "78÷6=", "26÷7="),
    @("32÷6=", "74÷4=", "72÷7=", "51÷2=", "29÷5=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowValues = $values[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $rowValues[$c - 1]
    }
}
